# This workbook holds one weekly record per row for "Cilantro" price data.
# The edit adds a new weekly record, inserted as row 83, which pushes every
# existing row from 83 downward down by one (old row 83 becomes row 84, ...,
# old row 216 becomes row 217). The sheet's used range grows from
# A1:R216 to A1:R217.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 83 (shifts rows 83:216 down to 84:217).
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new weekly record.
$ws.Range("A83").Value = 3
$ws.Range("B83").Value = "Femacal de La Calera"
$ws.Range("C83").Value = "Coquimbo"
$ws.Range("D83").Value = 44477
$ws.Range("E83").Value = 5
$ws.Range("F83").Value = 100112040
$ws.Range("G83").Value = "Cilantro"
$ws.Range("H83").Value = "Sin especificar"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 160
$ws.Range("K83").Value = 2500
$ws.Range("L83").Value = 2500
$ws.Range("M83").Value = 2500
$ws.Range("N83").Value = "$/docena de atados (3 kilos)"
$ws.Range("O83").Value = "Provincia de Quillota"
$ws.Range("P83").Value = 833
$ws.Range("Q83").Value = 3
$ws.Range("R83").Value = "Hortaliza"

# Give the date cell the same date number-format used by the rest of
# column D (style index 2 in the original workbook).
$ws.Range("D83").NumberFormat = $ws.Range("D84").NumberFormat()
